# Weekly price-sheet update: a new observation is inserted as row 28
# (pushing the existing rows 28-38 down to 29-39) and populated with the
# latest "Madrigal" Alcachofa price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 28, shifting rows 28-38
# down to 29-39 (matching the dimension growing from R38 to R39).
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly record.
$ws.Cells.Item(28, 1).Value  = 1
$ws.Cells.Item(28, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value  = 45093
$ws.Cells.Item(28, 5).Value  = 15
$ws.Cells.Item(28, 6).Value  = 100112013
$ws.Cells.Item(28, 7).Value  = "Alcachofa"
$ws.Cells.Item(28, 8).Value  = "Madrigal"
$ws.Cells.Item(28, 9).Value  = "Primera"
$ws.Cells.Item(28, 10).Value = 140
$ws.Cells.Item(28, 11).Value = 20000
$ws.Cells.Item(28, 12).Value = 22000
$ws.Cells.Item(28, 13).Value = 21000
$ws.Cells.Item(28, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(28, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(28, 16).Value = 525
$ws.Cells.Item(28, 17).Value = 40
$ws.Cells.Item(28, 18).Value = "Hortaliza"
